$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 1
    3  = 0
    4  = 1
    5  = 3
    6  = 0
    7  = 1
    8  = 0
    9  = 0
    10 = 0
    11 = 1
    12 = 0
    13 = 0
    14 = 0
    15 = 1
    16 = 0
    17 = 1
    18 = 0
    19 = 0
    20 = 1
    21 = 2
    23 = 1
    24 = 0
    25 = 1
    26 = 2
    27 = 2
    28 = 2
    29 = 0
    30 = 0
    31 = 0
    32 = 2
    33 = 1
    34 = 0
    35 = 0
    36 = 1
    37 = 1
    38 = 1
    39 = 0
    40 = 1
    41 = 2
    42 = 0
    43 = 0
    44 = 0
    45 = 2
    46 = 0
    47 = 0
    48 = 1
    49 = 0
    50 = 0
    51 = 3
    52 = 0
    53 = 1
    55 = 1
    56 = 2
    57 = 0
    58 = 0
    59 = 1
    60 = 1
    61 = 1
    62 = 0
    63 = 1
    64 = 1
    65 = 2
    66 = 1
    67 = 1
    68 = 2
    70 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
